$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 34; this shifts existing rows 34-68 down to 35-69
$ws.Rows.Item(34).Insert()

# Populate the newly inserted row 34 with the new data record
$ws.Range("A34").Value = 9
$ws.Range("B34").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C34").Value = "Metropolitana"
$ws.Range("D34").Value2 = 44484
$ws.Range("E34").Value = 13
$ws.Range("F34").Value = 100112022
$ws.Range("G34").Value = "Arveja Verde"
$ws.Range("H34").Value = "Perfection"
$ws.Range("I34").Value = "Primera"
$ws.Range("J34").Value = 43
$ws.Range("K34").Value = 22000
$ws.Range("L34").Value = 25000
$ws.Range("M34").Value = 23465
$ws.Range("N34").Value = "$/malla 25 kilos"
$ws.Range("O34").Value = "Provincia de Huasco"
$ws.Range("P34").Value = 939
$ws.Range("Q34").Value = 25
$ws.Range("R34").Value = "Hortaliza"
